$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two question rows (rows 6 and 7) - "Laravel" and "Recycler View" questions
$ws.Rows("6:7").Delete()

# --- Header row (row 2) ---
$ws.Range("E2").Value = "Bobot Pilgan"
$ws.Range("F2").Value = "Bobot Essay Praktikum"

# --- Row 3 : multiple-choice question ---
$ws.Range("C3").Value = "Pilgan"
$ws.Range("D3").Value = "Supervised Learning,Unsupervised Learning,Reinforcement Learning, Deep Learning"
$ws.Range("E3").Value = "5,0,0,0"
$ws.Range("F3").ClearContents()

# --- Row 4 : essay question ---
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = 15

# --- Row 5 : now the "Praktikum" question (replaces old JVM essay question) ---
$ws.Range("B5").Value = "Buatkan program java dengan tema apotek!"
$ws.Range("C5").Value = "Praktikum"
$ws.Range("E5").ClearContents()
$ws.Range("F5").Value = 10

# --- sheet view adjustments: scroll so column D is the left-most visible
#     column, and leave the active selection on F6 ---
$excel.Goto($ws.Range("D1"), $true)
$ws.Range("F6").Select()

# --- column width adjustment (column E) ---
$ws.Columns("E:E").ColumnWidth = 11.67
